$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - F column "想去人数" (want-to-go count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 98
$ws1.Range("F4").Value = 51
$ws1.Range("F6").Value = 143
$ws1.Range("F7").Value = 351
$ws1.Range("F8").Value = 5037
$ws1.Range("F10").Value = 5265
$ws1.Range("F11").Value = 606
$ws1.Range("F12").Value = 1330
$ws1.Range("F13").Value = 99

# Sheet "全部类型" (sheet4) - same F column updates (one extra row offset by +1 vs sheet1)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 98
$ws4.Range("F4").Value = 51
$ws4.Range("F6").Value = 143
$ws4.Range("F8").Value = 351
$ws4.Range("F9").Value = 5037
$ws4.Range("F11").Value = 5265
$ws4.Range("F12").Value = 606
$ws4.Range("F13").Value = 1330
$ws4.Range("F14").Value = 99

$wb.Save()
